# "added new search engine feature"
#
# A new record is inserted at row 3 (the result surfaced by the new search
# feature), and the record that used to live in row 3 ("Kontrak Asana" / AHM
# contract) is pushed down to row 4. F4 keeps its original value
# (C:/Users/Hp/Pictures/3 celah.pdf) — it is not part of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to stay a text cell even when its content looks like a
    # plain number (e.g. "1", "2"), matching the workbook's inline-string
    # storage instead of letting COM auto-coerce it to a numeric value.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 4 <- old row 3 content ("Kontrak Asana" contract).
$ws.Range("A4").Value = "Kontrak Asana"
$ws.Range("B4").Value = "AHM"
Set-TextValue "C4" "1"
$ws.Range("D4").Value = "16-06-2020 21:00:00"
$ws.Range("E4").Value = "16-06-2020 22:46:00"
# F4 unchanged: "C:/Users/Hp/Pictures/3 celah.pdf"

# Row 3 <- new record added by the search engine feature.
Set-TextValue "A3" "1"
$ws.Range("B3").Value = "abadi ehe"
Set-TextValue "C3" "2"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "E:/16719147 Arsyi Adlani Introductory Paragraph.pdf"
